$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Refresh the "last_edited_time" (column D) timestamps for this export.
$ws.Range("D2").Value = "2024-08-31T05:43:00.000Z"
$ws.Range("D3:D22").Value = "2024-08-31T05:40:00.000Z"
